$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 75, shifting existing rows 75-153 down to 76-154.
$ws.Rows.Item(75).Insert()

# Fill in the new row 75 with the new weekly data point.
$ws.Range("A75").Value = 4
$ws.Range("B75").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C75").Value = "Los Lagos"
$ws.Range("D75").Value = 44907
$ws.Range("E75").Value = 10
$ws.Range("F75").Value = 100112052
$ws.Range("G75").Value = "Albahaca"
$ws.Range("H75").Value = "Sin especificar"
$ws.Range("I75").Value = "Primera"
$ws.Range("J75").Value = 40
$ws.Range("K75").Value = 8000
$ws.Range("L75").Value = 8000
$ws.Range("M75").Value = 8000
$ws.Range("N75").Value = "$/docena de matas"
$ws.Range("O75").Value = "Región Metropolitana"
$ws.Range("P75").Value = 1333
$ws.Range("Q75").Value = 6
$ws.Range("R75").Value = "Hortaliza"

# D column stores dates, keep its date-like number format consistent with the rest of column D.
$ws.Range("D75").NumberFormat = $ws.Range("D76").NumberFormat
